$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Lgi4"
$ws.Cells.Item(2,3).Value = "Adam22"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.06754433333333333
$ws.Cells.Item(2,8).Value = 0.202633
$ws.Cells.Item(2,9).Value = 0.02266275549884949
$ws.Cells.Item(2,10).Value = 0.02266275549884949
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.357976666666667
$ws.Cells.Item(2,14).Value = 7.07393
$ws.Cells.Item(2,15).Value = 0.4138076229453391
$ws.Cells.Item(2,16).Value = 0.4138076229453391
$ws.Cells.Item(2,17).Value = 0.1592679619655555
$ws.Cells.Item(2,18).Value = 1.43341165769
$ws.Cells.Item(2,19).Value = 0.00937802098237032
$ws.Cells.Item(2,20).Value = 0.00937802098237032

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Lgi4"
$ws.Cells.Item(3,3).Value = "Adam22"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.06754433333333333
$ws.Cells.Item(3,8).Value = 0.202633
$ws.Cells.Item(3,9).Value = 0.02266275549884949
$ws.Cells.Item(3,10).Value = 0.02266275549884949
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.653161666666667
$ws.Cells.Item(3,14).Value = 4.959485
$ws.Cells.Item(3,15).Value = 0.2901177561670903
$ws.Cells.Item(3,16).Value = 0.2901177561670903
$ws.Cells.Item(3,17).Value = 0.1116617026672222
$ws.Cells.Item(3,18).Value = 1.004955324005
$ws.Cells.Item(3,19).Value = 0.006574867773889601
$ws.Cells.Item(3,20).Value = 0.006574867773889602

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Lgi4"
$ws.Cells.Item(4,3).Value = "Adam22"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.06754433333333333
$ws.Cells.Item(4,8).Value = 0.202633
$ws.Cells.Item(4,9).Value = 0.02266275549884949
$ws.Cells.Item(4,10).Value = 0.02266275549884949
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 0.33798
$ws.Cells.Item(4,14).Value = 1.01394
$ws.Cells.Item(4,15).Value = 0.05931301288098655
$ws.Cells.Item(4,16).Value = 0.05931301288098655
$ws.Cells.Item(4,17).Value = 0.02282863378
$ws.Cells.Item(4,18).Value = 0.20545770402
$ws.Cells.Item(4,19).Value = 0.001344196308821908
$ws.Cells.Item(4,20).Value = 0.001344196308821909

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Lgi4"
$ws.Cells.Item(5,3).Value = "Adam22"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.06754433333333333
$ws.Cells.Item(5,8).Value = 0.202633
$ws.Cells.Item(5,9).Value = 0.02266275549884949
$ws.Cells.Item(5,10).Value = 0.02266275549884949
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.5223656666666666
$ws.Cells.Item(5,14).Value = 1.567097
$ws.Cells.Item(5,15).Value = 0.091671345983742
$ws.Cells.Item(5,16).Value = 0.091671345983742
$ws.Cells.Item(5,17).Value = 0.03528284071122222
$ws.Cells.Item(5,18).Value = 0.317545566401
$ws.Cells.Item(5,19).Value = 0.002077525300279983
$ws.Cells.Item(5,20).Value = 0.002077525300279983

# Row 6
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Lgi4"
$ws.Cells.Item(6,3).Value = "Adam22"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.06754433333333333
$ws.Cells.Item(6,8).Value = 0.202633
$ws.Cells.Item(6,9).Value = 0.02266275549884949
$ws.Cells.Item(6,10).Value = 0.02266275549884949
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.8267596666666667
$ws.Cells.Item(6,14).Value = 2.480279
$ws.Cells.Item(6,15).Value = 0.145090262022842
$ws.Cells.Item(6,16).Value = 0.145090262022842
$ws.Cells.Item(6,17).Value = 0.05584293051188888
$ws.Cells.Item(6,18).Value = 0.502586374607
$ws.Cells.Item(6,19).Value = 0.003288145133487675
$ws.Cells.Item(6,20).Value = 0.003288145133487676

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Lgi4"
$ws.Cells.Item(7,3).Value = "Adam22"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.302381666666667
$ws.Cells.Item(7,8).Value = 6.907145
$ws.Cells.Item(7,9).Value = 0.7725046677002302
$ws.Cells.Item(7,10).Value = 0.7725046677002302
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.357976666666667
$ws.Cells.Item(7,14).Value = 7.07393
$ws.Cells.Item(7,15).Value = 0.4138076229453391
$ws.Cells.Item(7,16).Value = 0.4138076229453391
$ws.Cells.Item(7,17).Value = 5.428962247761111
$ws.Cells.Item(7,18).Value = 48.86066022985
$ws.Cells.Item(7,19).Value = 0.3196683202552114
$ws.Cells.Item(7,20).Value = 0.3196683202552114

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Lgi4"
$ws.Cells.Item(8,3).Value = "Adam22"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.302381666666667
$ws.Cells.Item(8,8).Value = 6.907145
$ws.Cells.Item(8,9).Value = 0.7725046677002302
$ws.Cells.Item(8,10).Value = 0.7725046677002302
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 1.653161666666667
$ws.Cells.Item(8,14).Value = 4.959485
$ws.Cells.Item(8,15).Value = 0.2901177561670903
$ws.Cells.Item(8,16).Value = 0.2901177561670903
$ws.Cells.Item(8,17).Value = 3.806209113369444
$ws.Cells.Item(8,18).Value = 34.255882020325
$ws.Cells.Item(8,19).Value = 0.2241173208217945
$ws.Cells.Item(8,20).Value = 0.2241173208217945

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Lgi4"
$ws.Cells.Item(9,3).Value = "Adam22"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.302381666666667
$ws.Cells.Item(9,8).Value = 6.907145
$ws.Cells.Item(9,9).Value = 0.7725046677002302
$ws.Cells.Item(9,10).Value = 0.7725046677002302
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.33798
$ws.Cells.Item(9,14).Value = 1.01394
$ws.Cells.Item(9,15).Value = 0.05931301288098655
$ws.Cells.Item(9,16).Value = 0.05931301288098655
$ws.Cells.Item(9,17).Value = 0.7781589557
$ws.Cells.Item(9,18).Value = 7.003430601300001
$ws.Cells.Item(9,19).Value = 0.04581957930592599
$ws.Cells.Item(9,20).Value = 0.04581957930592599

# Row 10
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Lgi4"
$ws.Cells.Item(10,3).Value = "Adam22"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.302381666666667
$ws.Cells.Item(10,8).Value = 6.907145
$ws.Cells.Item(10,9).Value = 0.7725046677002302
$ws.Cells.Item(10,10).Value = 0.7725046677002302
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.5223656666666666
$ws.Cells.Item(10,14).Value = 1.567097
$ws.Cells.Item(10,15).Value = 0.091671345983742
$ws.Cells.Item(10,16).Value = 0.091671345983742
$ws.Cells.Item(10,17).Value = 1.202685134229444
$ws.Cells.Item(10,18).Value = 10.824166208065
$ws.Cells.Item(10,19).Value = 0.07081654266680346
$ws.Cells.Item(10,20).Value = 0.07081654266680346

# Row 11
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Lgi4"
$ws.Cells.Item(11,3).Value = "Adam22"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 2.302381666666667
$ws.Cells.Item(11,8).Value = 6.907145
$ws.Cells.Item(11,9).Value = 0.7725046677002302
$ws.Cells.Item(11,10).Value = 0.7725046677002302
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.8267596666666667
$ws.Cells.Item(11,14).Value = 2.480279
$ws.Cells.Item(11,15).Value = 0.145090262022842
$ws.Cells.Item(11,16).Value = 0.145090262022842
$ws.Cells.Item(11,17).Value = 1.903516299272778
$ws.Cells.Item(11,18).Value = 17.131646693455
$ws.Cells.Item(11,19).Value = 0.1120829046504949
$ws.Cells.Item(11,20).Value = 0.1120829046504949

# Row 12
$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Lgi4"
$ws.Cells.Item(12,3).Value = "Adam22"
$ws.Cells.Item(12,4).Value = "ECs"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.6104853333333333
$ws.Cells.Item(12,8).Value = 1.831456
$ws.Cells.Item(12,9).Value = 0.2048325768009203
$ws.Cells.Item(12,10).Value = 0.2048325768009203
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.357976666666667
$ws.Cells.Item(12,14).Value = 7.07393
$ws.Cells.Item(12,15).Value = 0.4138076229453391
$ws.Cells.Item(12,16).Value = 0.4138076229453391
$ws.Cells.Item(12,17).Value = 1.439510171342222
$ws.Cells.Item(12,18).Value = 12.95559154208
$ws.Cells.Item(12,19).Value = 0.08476128170775744
$ws.Cells.Item(12,20).Value = 0.08476128170775746

# Row 13
$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Lgi4"
$ws.Cells.Item(13,3).Value = "Adam22"
$ws.Cells.Item(13,4).Value = "FAPs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.6104853333333333
$ws.Cells.Item(13,8).Value = 1.831456
$ws.Cells.Item(13,9).Value = 0.2048325768009203
$ws.Cells.Item(13,10).Value = 0.2048325768009203
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.653161666666667
$ws.Cells.Item(13,14).Value = 4.959485
$ws.Cells.Item(13,15).Value = 0.2901177561670903
$ws.Cells.Item(13,16).Value = 0.2901177561670903
$ws.Cells.Item(13,17).Value = 1.009230951128889
$ws.Cells.Item(13,18).Value = 9.08307856016
$ws.Cells.Item(13,19).Value = 0.0594255675714062
$ws.Cells.Item(13,20).Value = 0.05942556757140621

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Lgi4"
$ws.Cells.Item(14,3).Value = "Adam22"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.6104853333333333
$ws.Cells.Item(14,8).Value = 1.831456
$ws.Cells.Item(14,9).Value = 0.2048325768009203
$ws.Cells.Item(14,10).Value = 0.2048325768009203
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 0.33798
$ws.Cells.Item(14,14).Value = 1.01394
$ws.Cells.Item(14,15).Value = 0.05931301288098655
$ws.Cells.Item(14,16).Value = 0.05931301288098655
$ws.Cells.Item(14,17).Value = 0.20633183296
$ws.Cells.Item(14,18).Value = 1.85698649664
$ws.Cells.Item(14,19).Value = 0.01214923726623865
$ws.Cells.Item(14,20).Value = 0.01214923726623865

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Lgi4"
$ws.Cells.Item(15,3).Value = "Adam22"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.6104853333333333
$ws.Cells.Item(15,8).Value = 1.831456
$ws.Cells.Item(15,9).Value = 0.2048325768009203
$ws.Cells.Item(15,10).Value = 0.2048325768009203
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 0.5223656666666666
$ws.Cells.Item(15,14).Value = 1.567097
$ws.Cells.Item(15,15).Value = 0.091671345983742
$ws.Cells.Item(15,16).Value = 0.091671345983742
$ws.Cells.Item(15,17).Value = 0.3188965781368889
$ws.Cells.Item(15,18).Value = 2.870069203232
$ws.Cells.Item(15,19).Value = 0.01877727801665857
$ws.Cells.Item(15,20).Value = 0.01877727801665857

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Lgi4"
$ws.Cells.Item(16,3).Value = "Adam22"
$ws.Cells.Item(16,4).Value = "sCs"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.6104853333333333
$ws.Cells.Item(16,8).Value = 1.831456
$ws.Cells.Item(16,9).Value = 0.2048325768009203
$ws.Cells.Item(16,10).Value = 0.2048325768009203
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.8267596666666667
$ws.Cells.Item(16,14).Value = 2.480279
$ws.Cells.Item(16,15).Value = 0.145090262022842
$ws.Cells.Item(16,16).Value = 0.145090262022842
$ws.Cells.Item(16,17).Value = 0.5047246506915556
$ws.Cells.Item(16,18).Value = 4.542521856224
$ws.Cells.Item(16,19).Value = 0.02971921223885943
$ws.Cells.Item(16,20).Value = 0.02971921223885944
